# Add two new daily rows (116, 117 -> dates 2025-11-24 / 2025-11-25, serials
# 45985 / 45986) to every sheet in the workbook, mirroring the pattern
# already used for each prior row (date in column A, value in column B).

$wb = $excel.ActiveWorkbook

# New B-column values for row 116 and row 117, one pair per worksheet, in
# worksheet order (sheet1..sheet7 == Worksheets.Item(1..7)).
$newValues = @(
    @(3193, 3169),   # 진양산업
    @(1081, 1137),   # 넥스트아이
    @(1212, 1200),   # 삼보산업
    @(1845, 1993),   # YBM넷
    @(724,  722),    # NE능률
    @(1537, 1532),   # 위즈코프
    @(2577, 2559)    # 대영포장
)

$dateRow116 = 45985
$dateRow117 = 45986

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $pair = $newValues[$i - 1]

    $lastRow = 115

    # Column A: dates, carrying the same style as the preceding date cell.
    $ws.Cells.Item($lastRow + 1, 1).Value = $dateRow116
    $ws.Cells.Item($lastRow + 1, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

    $ws.Cells.Item($lastRow + 2, 1).Value = $dateRow117
    $ws.Cells.Item($lastRow + 2, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

    # Column B: the new remn_amt values.
    $ws.Cells.Item($lastRow + 1, 2).Value = $pair[0]
    $ws.Cells.Item($lastRow + 2, 2).Value = $pair[1]
}
